# Weekly Fruta/Hortaliza update: add a new price-report week for
# Femacal de La Calera - Repollo (Crespo record), both "Primera" and
# "Segunda" quality rows, dated 2021-11-09 (Excel serial 44509).
#
# This pushes the existing data block (old rows 375-389) down by two
# rows (to 377-391) - which Rows.Insert() does for us automatically,
# carrying along values/formatting - and then fills in the two newly
# opened rows (375-376) with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 375, shifting the rest
# of the table (old 375:389) down to 377:391.
$ws.Rows.Item(375).Insert()
$ws.Rows.Item(375).Insert()

# New row 375: Repollo, Crespo record, "Primera"
$ws.Cells.Item(375, 1).Value = 3
$ws.Cells.Item(375, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(375, 3).Value = "Coquimbo"
$ws.Cells.Item(375, 4).Value = [datetime]"2021-11-09"
$ws.Cells.Item(375, 5).Value = 5
$ws.Cells.Item(375, 6).Value = 100112006
$ws.Cells.Item(375, 7).Value = "Repollo"
$ws.Cells.Item(375, 8).Value = "Crespo record"
$ws.Cells.Item(375, 9).Value = "Primera"
$ws.Cells.Item(375, 10).Value = 880
$ws.Cells.Item(375, 11).Value = 600
$ws.Cells.Item(375, 12).Value = 600
$ws.Cells.Item(375, 13).Value = 600
$ws.Cells.Item(375, 14).Value = "`$/unidad"
$ws.Cells.Item(375, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(375, 16).Value = 600
$ws.Cells.Item(375, 17).Value = 1
$ws.Cells.Item(375, 18).Value = "Hortaliza"

# New row 376: Repollo, Crespo record, "Segunda"
$ws.Cells.Item(376, 1).Value = 3
$ws.Cells.Item(376, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(376, 3).Value = "Coquimbo"
$ws.Cells.Item(376, 4).Value = [datetime]"2021-11-09"
$ws.Cells.Item(376, 5).Value = 5
$ws.Cells.Item(376, 6).Value = 100112006
$ws.Cells.Item(376, 7).Value = "Repollo"
$ws.Cells.Item(376, 8).Value = "Crespo record"
$ws.Cells.Item(376, 9).Value = "Segunda"
$ws.Cells.Item(376, 10).Value = 900
$ws.Cells.Item(376, 11).Value = 500
$ws.Cells.Item(376, 12).Value = 500
$ws.Cells.Item(376, 13).Value = 500
$ws.Cells.Item(376, 14).Value = "`$/unidad"
$ws.Cells.Item(376, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(376, 16).Value = 500
$ws.Cells.Item(376, 17).Value = 1
$ws.Cells.Item(376, 18).Value = "Hortaliza"
